# Fix typos in the homework description (task numbers / lab number).
$d = $word.ActiveDocument

function Replace-ExactText($findText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        # Assign Range.Text directly (rather than using Find's own ReplaceWith)
        # so Word's smart-quote AutoCorrect doesn't mangle the literal apostrophe.
        $r.Text = $newText
    } else {
        Write-Host "WARNING: text not found -> $findText"
    }
}

# "Start with Task #3 ..." -> "Start with Task #2 ..."
Replace-ExactText "Start with Task #3 as it is by far the easiest" "Start with Task #2 as it is by far the easiest"

# "Next, try for Task #4A ..." -> "Next, try for Task #3 ..."
Replace-ExactText "Next, try for Task #4A as it's the most straight forward of the remaining tasks. " "Next, try for Task #3 as it's the most straight forward of the remaining tasks. "

# "Depending on your implementation of Task #4A, completing Task #4B ... looping back to 4B ..."
# -> "... Task #3, completing Task #4 ... looping back to 4 ..."
Replace-ExactText "Depending on your implementation of Task #4A, completing Task #4B may be rather difficult.  If this is the case, consider switching to Task #5 and looping back to 4B if you have the time. " "Depending on your implementation of Task #3, completing Task #4 may be rather difficult.  If this is the case, consider switching to Task #5 and looping back to 4 if you have the time. "

# "Task #4B and Task #5 are related ..." -> "Task #4 and Task #5 are related ..."
Replace-ExactText "Task #4B and Task #5 are related and have the potential to share a lot of code." "Task #4 and Task #5 are related and have the potential to share a lot of code."

# "During lab 11, you must demonstrate" -> "During lab 12, you must demonstrate"
Replace-ExactText "During lab 11, you must demonstrate" "During lab 12, you must demonstrate"
